$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: update timestamp (C9)
$ws.Range("C9").Value = 43380.823182870372

# Row 10: sensor name Sensor5 -> Sensor2, timestamp update
$ws.Range("B10").Value = "Sensor2"
$ws.Range("C10").Value = 43381.313101851854

# Row 11: sensor name Sensor5 -> Sensor2, timestamp update
$ws.Range("B11").Value = "Sensor2"
$ws.Range("C11").Value = 43381.313101851854

# Row 12: technology LunaSensor -> CarendoSensor, sensor Sensor5 -> Sensor1,
# timestamp update, bemaerkning Ekstra personale tilstede -> Borger kræver to personaler
$ws.Range("A12").Value = "CarendoSensor"
$ws.Range("B12").Value = "Sensor1"
$ws.Range("C12").Value = 43380.620393518519
$ws.Range("D12").Value = "Borger kræver to personaler"

# Update column widths (col A widened for "CarendoSensor", col D widened for the
# new longer "Borger kraever to personaler" bemaerkning; col E is left untouched so
# it keeps its original width/bestFit)
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 23.333333333333332

# Update selection to E14
$ws.Range("E14").Select()
